# Schematic Updates & Parts List Updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---
$ws.Name = "DSP Board"

# --- Column widths for the new parts-list table ---
$ws.Columns("A").ColumnWidth = 25.7109375
$ws.Columns("B").ColumnWidth = 20.7109375
$ws.Columns("C").ColumnWidth = 22.5703125
$ws.Columns("D").ColumnWidth = 17.7109375
$ws.Columns("E").ColumnWidth = 50.7109375
$ws.Columns("F").ColumnWidth = 17.85546875
$ws.Columns("G").ColumnWidth = 41.42578125

# --- Move the "**" footnote from Q8 to S8 ---
$ws.Range("Q8").Cut($ws.Range("S8"))

# --- Give B9's empty neighbour the Hyperlink look (matches the existing link cell) ---
$ws.Range("C9").Style = "Hyperlink"

# --- Formatting for the new parts-list block ---
$ws.Range("A18").HorizontalAlignment = -4152
$ws.Range("B18:H18").HorizontalAlignment = -4108
$ws.Range("B19:H19").HorizontalAlignment = -4108
$ws.Range("B20:H20").HorizontalAlignment = -4108
$ws.Range("B21:H21").HorizontalAlignment = -4108
$ws.Range("B22:E22").HorizontalAlignment = -4108
$ws.Range("G22:H22").HorizontalAlignment = -4108
$ws.Range("D23:E23").HorizontalAlignment = -4108
$ws.Range("G23:H23").HorizontalAlignment = -4108
$ws.Range("D24:E24").HorizontalAlignment = -4108
$ws.Range("G24:H24").HorizontalAlignment = -4108
$ws.Range("B25:H38").HorizontalAlignment = -4108

# "Bad" (red) cells for unsourced parts - style + centering applied together
$ws.Range("B23").Style = "Bad"
$ws.Range("B23").HorizontalAlignment = -4108
$ws.Range("C23").Style = "Bad"
$ws.Range("C23").HorizontalAlignment = -4108
$ws.Range("F23").Style = "Bad"
$ws.Range("F23").HorizontalAlignment = -4108
$ws.Range("B24").Style = "Bad"
$ws.Range("B24").HorizontalAlignment = -4108
$ws.Range("C24").Style = "Bad"
$ws.Range("C24").HorizontalAlignment = -4108
$ws.Range("F24").Style = "Bad"
$ws.Range("F24").HorizontalAlignment = -4108

# --- Table text, written in the same order it was originally typed so that
#     the shared-string table indices line up with the source file ---
$ws.Range("A18").Value = "Item Number"
$ws.Range("B18").Value = "Item Name"
$ws.Range("D18").Value = "Quantity"
$ws.Range("E18").Value = "Purpose"
$ws.Range("G18").Value = "Source"
$ws.Range("H18").Value = "Price"
$ws.Range("C18").Value = "Footprint"
$ws.Range("F18").Value = "Manufacturer"
$ws.Range("F19").Value = "TI"
$ws.Range("C19").Value = "176-Pin QFP"
$ws.Range("B19").Value = "TMS320F28335"
$ws.Range("E19").Value = "Digital Signal Processor"
$ws.Range("B20").Value = "ADS8320"
$ws.Range("C20").Value = "MSOP-8"
$ws.Range("E20").Value = "Audio Sampler"
$ws.Range("B21").Value = "TPS70351"
$ws.Range("C21").Value = "24-Pin PowerPAD TSSOP"
$ws.Range("E21").Value = "5.0V to 3.3V/1.8V Power Delivery and Sequencing"
$ws.Range("B22").Value = "Epson Oscillator"
$ws.Range("C22").Value = "4-Pin (Unique)"
$ws.Range("E22").Value = "Oscillator Source for DSP Clock"
$ws.Range("E23").Value = "Audio BandPass Filter Op-Amps"
$ws.Range("E24").Value = "Audio PreAmp Op-Amps"

# --- F20 / F21 reuse the shared "TI" string already created above ---
$ws.Range("F20").Value = "TI"
$ws.Range("F21").Value = "TI"

# --- Quantity column (numbers, not shared strings) ---
$ws.Range("D19").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("D24").Value = 1

# --- Item Number column for the populated + placeholder rows ---
$ws.Range("A19").Value = 1
$ws.Range("A20").Value = 2
$ws.Range("A21").Value = 3
$ws.Range("A22").Value = 4
$ws.Range("A23").Value = 5
$ws.Range("A24").Value = 6
$ws.Range("A25").Value = 7
$ws.Range("A26").Value = 8
$ws.Range("A27").Value = 9
$ws.Range("A28").Value = 10
$ws.Range("A29").Value = 11
$ws.Range("A30").Value = 12
$ws.Range("A31").Value = 13
$ws.Range("A32").Value = 14
$ws.Range("A33").Value = 15
$ws.Range("A34").Value = 16
$ws.Range("A35").Value = 17
$ws.Range("A36").Value = 18
$ws.Range("A37").Value = 19
$ws.Range("A38").Value = 20

# --- Move the view to where the new table lives ---
$win = $excel.Windows.Item(1)
$win.ScrollRow = 4
$win.ScrollColumn = 1
$excel.Goto($ws.Range("A17"))
